# Regenerate merged AHB files:
#  - rename header row from the "_old" / "_new" suffix convention to the
#    "_FV2410" / "_FV2504" convention used by the new merge script
#  - turn the data range into a proper Excel Table ("Table1")
#  - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (A1:J1 = "before" block, K1 = diff marker,
#        L1:U1 = "after" block) -------------------------------------------
$ws.Range("A1").Value = "Segmentname_FV2410"
$ws.Range("B1").Value = "Segmentgruppe_FV2410"
$ws.Range("C1").Value = "Segment_FV2410"
$ws.Range("D1").Value = "Datenelement_FV2410"
$ws.Range("E1").Value = "Segment ID_FV2410"
$ws.Range("F1").Value = "Code_FV2410"
$ws.Range("G1").Value = "Qualifier_FV2410"
$ws.Range("H1").Value = "Beschreibung_FV2410"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("J1").Value = "Bedingung_FV2410"

$ws.Range("L1").Value = "Segmentname_FV2504"
$ws.Range("M1").Value = "Segmentgruppe_FV2504"
$ws.Range("N1").Value = "Segment_FV2504"
$ws.Range("O1").Value = "Datenelement_FV2504"
$ws.Range("P1").Value = "Segment ID_FV2504"
$ws.Range("Q1").Value = "Code_FV2504"
$ws.Range("R1").Value = "Qualifier_FV2504"
$ws.Range("S1").Value = "Beschreibung_FV2504"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2504"
$ws.Range("U1").Value = "Bedingung_FV2504"

# --- 2. Turn A1:U72 into a native Excel table named "Table1" --------------
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U72"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# --- 3. Freeze the header row ---------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
